# Updated Master data as per 16th May Refresh
# Append three new user rows (34-36) to Sheet1, mirroring the existing
# row layout (id, uin, name, email, mobile, status_code, lang_code,
# last_login_method, is_active, cr_by, cr_dtimes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new "name" values first (so shared-string order matches source data,
#     which groups all names together before all emails) ---
$ws.Range("C34").Value = "Nikola Tesla"
$ws.Range("C35").Value = "Graham Bell"
$ws.Range("C36").Value = "Albert Miles"

# --- then the "email" values ---
$ws.Range("D34").Value = "nikola.tesla@xyz.com"
$ws.Range("D35").Value = "graham.bell@xyz.com"
$ws.Range("D36").Value = "albert.miles@xyz.com"

# --- id / uin / mobile (plain numbers) ---
$ws.Range("A34").Value = 110033
$ws.Range("B34").Value = 9317596771
$ws.Range("E34").Value = 818876434

$ws.Range("A35").Value = 110034
$ws.Range("B35").Value = 9317596772
$ws.Range("E35").Value = 818876435

$ws.Range("A36").Value = 110035
$ws.Range("B36").Value = 9317596773
$ws.Range("E36").Value = 818876436

# --- status_code / lang_code / last_login_method / cr_by / cr_dtimes
#     (re-use existing shared strings, same for every row) ---
foreach ($r in 34..36) {
    $ws.Range("F$r").Value = "ACT"
    $ws.Range("G$r").Value = "eng"
    $ws.Range("H$r").Value = "PWD"
    $ws.Range("J$r").Value = "superadmin"
    $ws.Range("K$r").Value = "now()"
}

# --- is_active (boolean, left-aligned like the rest of the column) ---
$ws.Range("I34:I36").Value = $true
$ws.Range("I34:I36").HorizontalAlignment = -4131

# Match the post-edit selection left by Excel after the new rows were
# added (whole-row selection starting right below the new data).
[void]$ws.Range("A37:XFD1048576").Select()
